$d = $word.ActiveDocument

# 1) "Power-ups (get bigger, enemies move slower, etc.)" becomes the
#    "More stuff for Points system (item collection - coins = 50 pts.)" item.
$enDash = [char]0x2013
$newFirst = "More stuff for Points system (item collection " + $enDash + " coins = 50 pts.)"
$null = $d.Content.Find.Execute(
    "Power-ups (get bigger, enemies move slower, etc.)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    $newFirst, 2)

# 2) The old "More stuff for ... etc.)" item becomes "Sound effects".
$null = $d.Content.Find.Execute(
    "More stuff for Points system (item collection, finish level, lives bonus when you finish, etc.)",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Sound effects", 2)

# 3) The old "Sound effects" item (now the third list entry) becomes
#    "Print scene for game over without finishing level".
# There are now two paragraphs containing "Sound effects" text (the one we
# just created and the original one); restrict the search to start after
# the paragraph we just edited so we hit the original one.
$searchStart = $d.Paragraphs(4).Range.Start
$searchRange = $d.Range($searchStart, $d.Content.End)
$null = $searchRange.Find.Execute(
    "Sound effects",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Print scene for game over without finishing level", 2)
